# "add 19 test cases" - append 7 new test-case rows (41-47) to the tracker
# sheet, each naming one of the newly authored test cases.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRows = @(
    @(41, "ChangeLanguage"),
    @(42, "SubmitDebugLogs"),
    @(43, "InvalidSecurityKey"),
    @(44, "SignOutThroughSetting"),
    @(45, "ChangeLanguageInSigninPage"),
    @(46, "ExploreRoomInFirstPage"),
    @(47, "CreateAccountBacktoSignInPage")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $name = $row[1]
    $ws.Cells.Item($r, 1).Value = $r
    $ws.Cells.Item($r, 2).Value = $name
}

# Restore the view/selection state recorded after the edit.
$win = $excel.ActiveWindow
$win.ScrollRow = 21
$win.ScrollColumn = 1
$ws.Range("C40").Select()
